$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '52.075.18'
$ws.Cells.Item(2, 5).Value = '  +0.23%  '
$ws.Cells.Item(3, 4).Value = '2.845.14'
$ws.Cells.Item(3, 5).Value = '  +2.57%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = '''363.21'
$ws.Cells.Item(5, 5).Value = '  +6.79%  '
$ws.Cells.Item(6, 4).Value = '''113.07'
$ws.Cells.Item(6, 5).Value = '  -2.11%  '
$ws.Cells.Item(7, 5).Value = '  +4.57%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$ws.Cells.Item(9, 5).Value = '  +4.89%  '
$ws.Cells.Item(10, 4).Value = '''41.58'
$ws.Cells.Item(10, 5).Value = '  -0.08%  '
$ws.Cells.Item(11, 4).Value = '''0.0865'
$ws.Cells.Item(11, 5).Value = '  +0.69%  '
$ws.Cells.Item(12, 5).Value = '  +0.89%  '
$ws.Cells.Item(13, 5).Value = '  +1.16%  '
$ws.Cells.Item(14, 5).Value = '  +3.14%  '
$ws.Cells.Item(15, 4).Value = '3.290.42'
$ws.Cells.Item(15, 5).Value = '  +2.31%  '
$ws.Cells.Item(16, 4).Value = '2.842.75'
$ws.Cells.Item(16, 5).Value = '  +1.79%  '
$ws.Cells.Item(17, 4).Value = '''0.916'
$ws.Cells.Item(17, 5).Value = '  +4.59%  '
$ws.Cells.Item(18, 4).Value = '52.183.21'
$ws.Cells.Item(18, 5).Value = '  +0.76%  '
$ws.Cells.Item(19, 5).Value = '  +9.14%  '
$ws.Cells.Item(20, 5).Value = '  -0.19%  '
$ws.Cells.Item(21, 4).Value = '''13.64'
$ws.Cells.Item(21, 5).Value = '  +3.20%  '
$ws.Cells.Item(22, 5).Value = '  +2.89%  '
$ws.Cells.Item(23, 5).Value = '  +0.80%  '
$ws.Cells.Item(24, 4).Value = '''269.39'
$ws.Cells.Item(24, 5).Value = '  -2.42%  '
$ws.Cells.Item(25, 4).Value = '''2.87'
$ws.Cells.Item(25, 5).Value = '  +4.77%  '
$ws.Cells.Item(26, 4).Value = '''27.15'
$ws.Cells.Item(26, 5).Value = '  +2.12%  '
$ws.Cells.Item(27, 5).Value = '  +0.10%  '
$ws.Cells.Item(28, 5).Value = '  +2.36%  '
$ws.Cells.Item(29, 5).Value = '  +1.27%  '
$ws.Cells.Item(30, 4).Value = '''0.0489'
$ws.Cells.Item(30, 5).Value = '  +30.31%  '
$ws.Cells.Item(31, 4).Value = '''53.85'
$ws.Cells.Item(31, 5).Value = '  +7.41%  '
$ws.Cells.Item(32, 5).Value = '  -0.02%  '
$ws.Cells.Item(33, 4).Value = '''35.44'
$ws.Cells.Item(33, 5).Value = '  +2.46%  '
$ws.Cells.Item(34, 4).Value = '''5.89'
$ws.Cells.Item(34, 5).Value = '  +3.59%  '
$ws.Cells.Item(35, 4).Value = '''5.50'
$ws.Cells.Item(35, 5).Value = '  +11.98%  '
$ws.Cells.Item(36, 4).Value = '''0.0847'
$ws.Cells.Item(36, 5).Value = '  +3.05%  '
$ws.Cells.Item(37, 5).Value = '  -0.09%  '
$ws.Cells.Item(38, 5).Value = '  +2.40%  '
$ws.Cells.Item(39, 4).Value = '''2.07'
$ws.Cells.Item(39, 5).Value = '  -0.98%  '
$ws.Cells.Item(40, 4).Value = '''18.52'
$ws.Cells.Item(40, 5).Value = '  -1.45%  '
$ws.Cells.Item(41, 4).Value = '''23.86'
$ws.Cells.Item(41, 5).Value = '  +3.84%  '
$ws.Cells.Item(42, 5).Value = '  +1.78%  '
$ws.Cells.Item(43, 4).Value = '''2.54'
$ws.Cells.Item(43, 5).Value = '  -3.39%  '
$ws.Cells.Item(44, 4).Value = '''126.54'
$ws.Cells.Item(44, 5).Value = '  +0.80%  '
$ws.Cells.Item(45, 5).Value = '  -2.72%  '
$ws.Cells.Item(46, 4).Value = '''3.43'
$ws.Cells.Item(46, 5).Value = '  +3.95%  '
$ws.Cells.Item(47, 4).Value = '2.115.92'
$ws.Cells.Item(47, 5).Value = '  +2.04%  '
$ws.Cells.Item(49, 4).Value = '''0.990'
$ws.Cells.Item(49, 5).Value = '  +14.24%  '
$ws.Cells.Item(50, 4).Value = '''5.90'
$ws.Cells.Item(50, 5).Value = '  +7.05%  '
$ws.Cells.Item(51, 4).Value = '''62.03'
$ws.Cells.Item(51, 5).Value = '  +5.03%  '
